# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.081.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.624.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.76"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.625.16"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.134"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.93%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.63"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.809.96"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.626.22"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "367.51"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.41"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.81"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.66%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.95"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.750.64"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "571.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.88"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.127"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.52"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.47"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.11"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.367"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.31"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0328"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +11.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.32"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "155.10"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.67"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.83"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.69"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.27%  "
